$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMSIN")

# Insert a new row above the current last data row (18). Excel inherits
# the formatting of the row above into the freshly inserted cells, which
# is how A18 and C18:G18 (previously unformatted / default style) pick up
# the sheet's normal data style - matching the diff, where those cells
# gain an explicit style. The record that used to live in row 18 is
# pushed down to row 19 by the insert.
$ws.Range("A18:G18").Insert(-4121)

# Re-enter the original row-18 record into the newly inserted row, fixing
# up the tiny floating point drift on the "Run Time" timestamp. The
# leading apostrophe keeps the date-look-alike "Run Date" text a literal
# string instead of Excel auto-converting it to a date serial number.
$ws.Range("A18").Value = "'2023-04-07"
$ws.Range("B18").Value = 45023.68523743055
$ws.Range("C18").Value = "176fstadh"
$ws.Range("D18").Value = 33
$ws.Range("E18").Value = 24
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 1.05

# Row 19 currently still holds the record shifted down from the old row
# 18; overwrite it with the new registration-history entry.
$ws.Range("A19").Value = "'2023-04-18"
$ws.Range("B19").Value = 45034.51305593146
$ws.Range("C19").Value = "176aadhtrail"
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = 33
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1.48
